$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The data rows (2-5) are cyclically shifted up by one: the old row 2 values
# move down to row 5, and old rows 3,4,5 move up to 2,3,4 — for columns
# D (Fecha), L (Calidad), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado) and S (Precio $/Kg).

$ws.Range("D2").Value = 44995
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 5500
$ws.Range("O2").Value = 6000
$ws.Range("P2").Value = 5750
$ws.Range("S2").Value = 2875

$ws.Range("D3").Value = 45008
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 60
$ws.Range("N3").Value = 7000
$ws.Range("O3").Value = 7000
$ws.Range("P3").Value = 7000
$ws.Range("S3").Value = 3500

$ws.Range("D4").Value = 45008
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 6000
$ws.Range("O4").Value = 6000
$ws.Range("P4").Value = 6000
$ws.Range("S4").Value = 3000

$ws.Range("D5").Value = 44991
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 6000
$ws.Range("O5").Value = 6000
$ws.Range("P5").Value = 6000
$ws.Range("S5").Value = 3000
